# Update the "want to go" counts (column F) on the relevant sheets to
# reflect newly generated output data.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 16326
$ws1.Range("F16").Value = 429
$ws1.Range("F20").Value = 608
$ws1.Range("F24").Value = 1151
$ws1.Range("F26").Value = 21
$ws1.Range("F32").Value = 82
$ws1.Range("F36").Value = 364
$ws1.Range("F39").Value = 5675

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 78

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 16326
$ws4.Range("F16").Value = 429
$ws4.Range("F20").Value = 608
$ws4.Range("F24").Value = 1151
$ws4.Range("F26").Value = 21
$ws4.Range("F31").Value = 78
$ws4.Range("F34").Value = 82
$ws4.Range("F38").Value = 364
$ws4.Range("F41").Value = 5675
